function Set-CellText {
    param($Worksheet, $CellRef, $Text)
    $range = $Worksheet.Range($CellRef)
    # Force a literal-text write: without pre-formatting the cell as
    # Text, Excel's smart type inference would coerce numeric-looking
    # strings like "318.42" or "3.88%" into real numbers/percentages
    # (and pick up a % number format), which would silently change both
    # the stored value and the cell's displayed formatting.
    $range.NumberFormat = "@"
    $range.Value = $Text
    # Revert to the workbook's default style now that the value is
    # safely stored as text, so we don't leave a stray text-format
    # style applied to a cell that originally had none.
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin-ranking snapshot (price / 1h volume %) for the existing
# rows, plus the upward shift of rows 8-17 caused by GateToken (GT)
# moving up in rank -- each coin in that block takes over the row
# previously occupied by the coin ranked one place above it.
Set-CellText $ws 'D2' '318.42'
Set-CellText $ws 'E2' '3.88%'
Set-CellText $ws 'D3' '36.20'
Set-CellText $ws 'E3' '0.00%'
Set-CellText $ws 'D4' '5.117'
Set-CellText $ws 'E4' '1.53%'
Set-CellText $ws 'D5' '0.08071'
Set-CellText $ws 'E5' '3.33%'
Set-CellText $ws 'D6' '2.174'
Set-CellText $ws 'E6' '2.64%'
Set-CellText $ws 'D7' '8.054'
Set-CellText $ws 'E7' '1.68%'
Set-CellText $ws 'B8' 'GateToken'
Set-CellText $ws 'C8' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-CellText $ws 'D8' '4.140'
Set-CellText $ws 'E8' '2.00%'
Set-CellText $ws 'B9' 'MXToken'
Set-CellText $ws 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 'D9' '0.9288'
Set-CellText $ws 'E9' '0.56%'
Set-CellText $ws 'B10' 'LiechtensteinCryptoassetsExchange'
Set-CellText $ws 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-CellText $ws 'D10' '0.1007'
Set-CellText $ws 'E10' '5.10%'
Set-CellText $ws 'B11' 'WazirX'
Set-CellText $ws 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-CellText $ws 'D11' '0.1875'
Set-CellText $ws 'E11' '-0.32%'
Set-CellText $ws 'B12' 'MandalaExchangeToken'
Set-CellText $ws 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-CellText $ws 'D12' '0.09150'
Set-CellText $ws 'E12' '5.00%'
Set-CellText $ws 'B13' 'BitrueCoin'
Set-CellText $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-CellText $ws 'D13' '0.03592'
Set-CellText $ws 'E13' '2.63%'
Set-CellText $ws 'B14' 'BitMartToken'
Set-CellText $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-CellText $ws 'D14' '0.09950'
Set-CellText $ws 'E14' '0.17%'
Set-CellText $ws 'B15' 'BitForexToken'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-CellText $ws 'D15' '0.001435'
Set-CellText $ws 'E15' '0.35%'
Set-CellText $ws 'B16' 'TigerCash'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-CellText $ws 'D16' '0.005702'
Set-CellText $ws 'E16' '0.10%'
Set-CellText $ws 'B17' 'LEO'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText $ws 'D17' '3.458'
Set-CellText $ws 'E17' '0.01%'
Set-CellText $ws 'D18' '2.799'
Set-CellText $ws 'E18' '15.96%'
Set-CellText $ws 'D19' '0.3374'
Set-CellText $ws 'E19' '-1.23%'
Set-CellText $ws 'D20' '0.1321'
Set-CellText $ws 'E20' '-1.83%'
Set-CellText $ws 'D21' '5.067'
Set-CellText $ws 'E21' '6.43%'
Set-CellText $ws 'D22' '0.2204'
Set-CellText $ws 'E22' '-3.93%'
Set-CellText $ws 'D23' '0.04603'
Set-CellText $ws 'E23' '-0.32%'
Set-CellText $ws 'D24' '0.001238'
Set-CellText $ws 'E24' '0.58%'
Set-CellText $ws 'D25' '0.004756'
Set-CellText $ws 'E25' '-6.76%'
Set-CellText $ws 'D26' '0.0001299'
Set-CellText $ws 'E26' '-7.41%'
Set-CellText $ws 'D39' '0.01945'
Set-CellText $ws 'E39' '6.42%'
Set-CellText $ws 'D40' '0.04955'
Set-CellText $ws 'E40' '4.17%'
Set-CellText $ws 'D41' '0.007814'
Set-CellText $ws 'E41' '3.92%'
Set-CellText $ws 'E42' '-0.38%'
Set-CellText $ws 'D43' '0.007560'
Set-CellText $ws 'E43' '-2.37%'
Set-CellText $ws 'D44' '0.002095'
Set-CellText $ws 'E44' '-6.24%'
Set-CellText $ws 'D45' '0.01167'
Set-CellText $ws 'E45' '12.16%'
Set-CellText $ws 'D46' '0.00006270'
Set-CellText $ws 'E46' '1.53%'
Set-CellText $ws 'E47' '-0.26%'
Set-CellText $ws 'D48' '28.90'
Set-CellText $ws 'E48' '-27.41%'
Set-CellText $ws 'D49' '0.001902'
Set-CellText $ws 'E49' '-5.10%'
Set-CellText $ws 'D50' '0.00002098'
Set-CellText $ws 'E50' '-0.26%'
Set-CellText $ws 'D51' '0.0001999'
Set-CellText $ws 'E51' '-0.26%'
